$wb = $excel.ActiveWorkbook

# ===================== Sheet: III year cse =====================
$ws = $wb.Worksheets.Item("III year cse")

# Insert first new row at 93 (ILAMUGUNTHAN N); copy style from row below, then set values later
$ws.Rows.Item(93).Insert()
$ws.Range("E94:K94").Copy()
$ws.Range("E93:K93").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Insert second new row at 121 (ABBINAVU T); copy style from row below
$ws.Rows.Item(121).Insert()
$ws.Range("E122:K122").Copy()
$ws.Range("E121:K121").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Write final values for rows 5-132 (SI, USERNAME, USER ID, SECTION, RANK, SOLVED, PAGE LINK)
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "PRAVINKUMAR S"
$ws.Range("G5").Value = "spravinkumar9952"
$ws.Range("H5").Value = "III year CSE B"
$ws.Range("I5").Value = 3529
$ws.Range("J5").Value = "3/4"
$ws.Range("K5").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/142/"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = "SANJEEV CHANDRAN M"
$ws.Range("G6").Value = "sanjeevchandran"
$ws.Range("H6").Value = "III year CSE B"
$ws.Range("I6").Value = 6252
$ws.Range("J6").Value = "3/4"
$ws.Range("K6").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/251/"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = "RAMMPRASHATH K"
$ws.Range("G7").Value = "ramm2413"
$ws.Range("H7").Value = "III year CSE B"
$ws.Range("I7").Value = 7020
$ws.Range("J7").Value = "3/4"
$ws.Range("K7").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/281/"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "BALAMURUGAN K"
$ws.Range("G8").Value = "user4029ok"
$ws.Range("H8").Value = "III year CSE A"
$ws.Range("I8").Value = 9480
$ws.Range("J8").Value = "3/4"
$ws.Range("K8").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/380/"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = "SASIREKA S"
$ws.Range("G9").Value = "sasireka20cs120"
$ws.Range("H9").Value = "III year CSE C"
$ws.Range("I9").Value = 9483
$ws.Range("J9").Value = "3/4"
$ws.Range("K9").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/380/"
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = "RAM DEEPAK P"
$ws.Range("G10").Value = "ramdeepak33"
$ws.Range("H10").Value = "III year CSE B"
$ws.Range("I10").Value = 10731
$ws.Range("J10").Value = "2/4"
$ws.Range("K10").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/430/"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = "TAMILSELVI S"
$ws.Range("G11").Value = "tamilselvi_53"
$ws.Range("H11").Value = "III year CSE C"
$ws.Range("I11").Value = 11921
$ws.Range("J11").Value = "2/4"
$ws.Range("K11").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/477/"
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = "SHREE RAMANAA M"
$ws.Range("G12").Value = "shree_ramanaa"
$ws.Range("H12").Value = "III year CSE C"
$ws.Range("I12").Value = 12073
$ws.Range("J12").Value = "2/4"
$ws.Range("K12").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/483/"
$ws.Range("E13").Value = 9
$ws.Range("F13").Value = "KAASIPRASANTH A"
$ws.Range("G13").Value = "kaasiprasanth_a"
$ws.Range("H13").Value = "III year CSE A"
$ws.Range("I13").Value = 12628
$ws.Range("J13").Value = "2/4"
$ws.Range("K13").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/506/"
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = "KAARTHIKEYAN A R"
$ws.Range("G14").Value = "a_r_kaarthikeyan"
$ws.Range("H14").Value = "III year CSE A"
$ws.Range("I14").Value = 13973
$ws.Range("J14").Value = "1/4"
$ws.Range("K14").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/559/"
$ws.Range("E15").Value = 11
$ws.Range("F15").Value = "SHAKTHI KIRAN R"
$ws.Range("G15").Value = "shakthi_kiran"
$ws.Range("H15").Value = "III year CSE C"
$ws.Range("I15").Value = 14173
$ws.Range("J15").Value = "1/4"
$ws.Range("K15").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/567/"
$ws.Range("E16").Value = 12
$ws.Range("F16").Value = "AISHWARYA S"
$ws.Range("G16").Value = "aishwarya_selvamurugan"
$ws.Range("H16").Value = "III year CSE A"
$ws.Range("I16").Value = 14299
$ws.Range("J16").Value = "1/4"
$ws.Range("K16").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/572/"
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = "SOUNDARYA V N"
$ws.Range("G17").Value = "soundarya_v_n"
$ws.Range("H17").Value = "III year CSE C"
$ws.Range("I17").Value = 14606
$ws.Range("J17").Value = "1/4"
$ws.Range("K17").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/585/"
$ws.Range("E18").Value = 14
$ws.Range("F18").Value = "VANITHA A"
$ws.Range("G18").Value = "vanitha_a"
$ws.Range("H18").Value = "III year CSE C"
$ws.Range("I18").Value = 14671
$ws.Range("J18").Value = "1/4"
$ws.Range("K18").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/587/"
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = "NITHYA M"
$ws.Range("G19").Value = "nithyamohan906"
$ws.Range("H19").Value = "III year CSE B"
$ws.Range("I19").Value = 14706
$ws.Range("J19").Value = "1/4"
$ws.Range("K19").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/589/"
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = "SELVANANDHINI A"
$ws.Range("G20").Value = "selvanandhini_a"
$ws.Range("H20").Value = "III year CSE C"
$ws.Range("I20").Value = 14724
$ws.Range("J20").Value = "1/4"
$ws.Range("K20").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/589/"
$ws.Range("E21").Value = 17
$ws.Range("F21").Value = "MADHAN KUMAR N"
$ws.Range("G21").Value = "madhankumar01"
$ws.Range("H21").Value = "III year CSE B"
$ws.Range("I21").Value = 14820
$ws.Range("J21").Value = "1/4"
$ws.Range("K21").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/593/"
$ws.Range("E22").Value = 18
$ws.Range("F22").Value = "PRAVIN D"
$ws.Range("G22").Value = "pravin-123"
$ws.Range("H22").Value = "III year CSE B"
$ws.Range("I22").Value = 14831
$ws.Range("J22").Value = "1/4"
$ws.Range("K22").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/594/"
$ws.Range("E23").Value = 19
$ws.Range("F23").Value = "RIHANABANU A"
$ws.Range("G23").Value = "rihanabanu"
$ws.Range("H23").Value = "III year CSE B"
$ws.Range("I23").Value = 14842
$ws.Range("J23").Value = "1/4"
$ws.Range("K23").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/594/"
$ws.Range("E24").Value = 20
$ws.Range("F24").Value = "PRAVEEN A"
$ws.Range("G24").Value = "praveen_a_"
$ws.Range("H24").Value = "III year CSE B"
$ws.Range("I24").Value = 14843
$ws.Range("J24").Value = "1/4"
$ws.Range("K24").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/594/"
$ws.Range("E25").Value = 21
$ws.Range("F25").Value = "PRIYA K"
$ws.Range("G25").Value = "priyakamaraj"
$ws.Range("H25").Value = "III year CSE B"
$ws.Range("I25").Value = 14852
$ws.Range("J25").Value = "1/4"
$ws.Range("K25").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/595/"
$ws.Range("E26").Value = 22
$ws.Range("F26").Value = "SANJITH R K"
$ws.Range("G26").Value = "sanjith16"
$ws.Range("H26").Value = "III year CSE B"
$ws.Range("I26").Value = 14853
$ws.Range("J26").Value = "1/4"
$ws.Range("K26").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/595/"
$ws.Range("E27").Value = 23
$ws.Range("F27").Value = "SOORJI MARTINA K"
$ws.Range("G27").Value = "soorji"
$ws.Range("H27").Value = "III year CSE C"
$ws.Range("I27").Value = 14860
$ws.Range("J27").Value = "1/4"
$ws.Range("K27").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/595/"
$ws.Range("E28").Value = 24
$ws.Range("F28").Value = "IJJU HEMANTH KUMAR"
$ws.Range("G28").Value = "ijjuhemanthkumar"
$ws.Range("H28").Value = "III year CSE A"
$ws.Range("I28").Value = 14873
$ws.Range("J28").Value = "1/4"
$ws.Range("K28").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/595/"
$ws.Range("E29").Value = 25
$ws.Range("F29").Value = "PRADEEP M"
$ws.Range("G29").Value = "pradeep_m_77"
$ws.Range("H29").Value = "III year CSE B"
$ws.Range("I29").Value = 14880
$ws.Range("J29").Value = "1/4"
$ws.Range("K29").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/596/"
$ws.Range("E30").Value = 26
$ws.Range("F30").Value = "RUBIN KUMAR K"
$ws.Range("G30").Value = "user9710zj"
$ws.Range("H30").Value = "III year CSE B"
$ws.Range("I30").Value = 14887
$ws.Range("J30").Value = "1/4"
$ws.Range("K30").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/596/"
$ws.Range("E31").Value = 27
$ws.Range("F31").Value = "SNEHA S"
$ws.Range("G31").Value = "user6961v"
$ws.Range("H31").Value = "III year CSE C"
$ws.Range("I31").Value = 14901
$ws.Range("J31").Value = "1/4"
$ws.Range("K31").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/597/"
$ws.Range("E32").Value = 28
$ws.Range("F32").Value = "SAKTHI S"
$ws.Range("G32").Value = "sakthi_s"
$ws.Range("H32").Value = "III year CSE B"
$ws.Range("I32").Value = 14903
$ws.Range("J32").Value = "1/4"
$ws.Range("K32").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/597/"
$ws.Range("E33").Value = 29
$ws.Range("F33").Value = "PAVITHRADEVI B"
$ws.Range("G33").Value = "pavithra_devi7"
$ws.Range("H33").Value = "III year CSE B"
$ws.Range("I33").Value = 14920
$ws.Range("J33").Value = "1/4"
$ws.Range("K33").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/597/"
$ws.Range("E34").Value = 30
$ws.Range("F34").Value = "INDHUMATHI B"
$ws.Range("G34").Value = "indhumathi_b"
$ws.Range("H34").Value = "III year CSE A"
$ws.Range("I34").Value = 14954
$ws.Range("J34").Value = "1/4"
$ws.Range("K34").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/599/"
$ws.Range("E35").Value = 31
$ws.Range("F35").Value = "PRAVEEN R"
$ws.Range("G35").Value = "r_praveen"
$ws.Range("H35").Value = "III year CSE B"
$ws.Range("I35").Value = 14979
$ws.Range("J35").Value = "1/4"
$ws.Range("K35").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/600/"
$ws.Range("E36").Value = 32
$ws.Range("F36").Value = "KISHORE P"
$ws.Range("G36").Value = "p_kishore"
$ws.Range("H36").Value = "III year CSE B"
$ws.Range("I36").Value = 16403
$ws.Range("J36").Value = "0/4"
$ws.Range("K36").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/657/"
$ws.Range("E37").Value = 33
$ws.Range("F37").Value = "SANJAY G S"
$ws.Range("G37").Value = "sanjudddd400"
$ws.Range("H37").Value = "III year CSE B"
$ws.Range("I37").Value = 18272
$ws.Range("J37").Value = "0/4"
$ws.Range("K37").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/731/"
$ws.Range("E38").Value = 34
$ws.Range("F38").Value = "KAVIN S"
$ws.Range("G38").Value = "kavin_s"
$ws.Range("H38").Value = "III year CSE A"
$ws.Range("I38").Value = 19233
$ws.Range("J38").Value = "0/4"
$ws.Range("K38").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/770/"
$ws.Range("E39").Value = 35
$ws.Range("F39").Value = "DEVAMITRA T"
$ws.Range("G39").Value = "devamitra_t"
$ws.Range("H39").Value = "III year CSE A"
$ws.Range("I39").Value = 19462
$ws.Range("J39").Value = "0/4"
$ws.Range("K39").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/779/"
$ws.Range("E40").Value = 36
$ws.Range("F40").Value = "INDIRA KUMAR A"
$ws.Range("G40").Value = "indiranj0"
$ws.Range("H40").Value = "III year CSE A"
$ws.Range("I40").Value = 19472
$ws.Range("J40").Value = "0/4"
$ws.Range("K40").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/779/"
$ws.Range("E41").Value = 37
$ws.Range("F41").Value = "SABAREESWARAN G"
$ws.Range("G41").Value = "sabari2309"
$ws.Range("H41").Value = "III year CSE B"
$ws.Range("I41").Value = 19473
$ws.Range("J41").Value = "0/4"
$ws.Range("K41").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/779/"
$ws.Range("E42").Value = 38
$ws.Range("F42").Value = "MENAGA E"
$ws.Range("G42").Value = "20cs071"
$ws.Range("H42").Value = "III year CSE B"
$ws.Range("I42").Value = 19487
$ws.Range("J42").Value = "0/4"
$ws.Range("K42").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/780/"
$ws.Range("E43").Value = 39
$ws.Range("F43").Value = "HONIKA S"
$ws.Range("G43").Value = "honi_ka_02"
$ws.Range("H43").Value = "III year CSE A"
$ws.Range("I43").Value = 19502
$ws.Range("J43").Value = "0/4"
$ws.Range("K43").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/781/"
$ws.Range("E44").Value = 40
$ws.Range("F44").Value = "YASWANTH KUMAR S"
$ws.Range("G44").Value = "user0208nr"
$ws.Range("H44").Value = "III year CSE C"
$ws.Range("I44").Value = 19512
$ws.Range("J44").Value = "0/4"
$ws.Range("K44").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/781/"
$ws.Range("E45").Value = 41
$ws.Range("F45").Value = "RAMASUBRAMANYAM P"
$ws.Range("G45").Value = "ram1216"
$ws.Range("H45").Value = "III year CSE B"
$ws.Range("I45").Value = 19521
$ws.Range("J45").Value = "0/4"
$ws.Range("K45").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/781/"
$ws.Range("E46").Value = 42
$ws.Range("F46").Value = "KAVIN R"
$ws.Range("G46").Value = "kavin-r"
$ws.Range("H46").Value = "III year CSE A"
$ws.Range("I46").Value = 19526
$ws.Range("J46").Value = "0/4"
$ws.Range("K46").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/782/"
$ws.Range("E47").Value = 43
$ws.Range("F47").Value = "PRATHIKSHA J"
$ws.Range("G47").Value = "prathi_03"
$ws.Range("H47").Value = "III year CSE B"
$ws.Range("I47").Value = 19561
$ws.Range("J47").Value = "0/4"
$ws.Range("K47").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/783/"
$ws.Range("E48").Value = 44
$ws.Range("F48").Value = "TAMILSELVAN M"
$ws.Range("G48").Value = "tamilselvan_2003"
$ws.Range("H48").Value = "III year CSE C"
$ws.Range("I48").Value = 19596
$ws.Range("J48").Value = "0/4"
$ws.Range("K48").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/784/"
$ws.Range("E49").Value = 45
$ws.Range("F49").Value = "GNANA SEKAR R"
$ws.Range("G49").Value = "gnanasekar_r"
$ws.Range("H49").Value = "III year CSE A"
$ws.Range("I49").Value = 19622
$ws.Range("J49").Value = "0/4"
$ws.Range("K49").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/785/"
$ws.Range("E50").Value = 46
$ws.Range("F50").Value = "MANORANJAN K"
$ws.Range("G50").Value = "manoranjank24_-"
$ws.Range("H50").Value = "III year CSE B"
$ws.Range("I50").Value = 19640
$ws.Range("J50").Value = "0/4"
$ws.Range("K50").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/786/"
$ws.Range("E51").Value = 47
$ws.Range("F51").Value = "PRAKASH RAJ S"
$ws.Range("G51").Value = "prakashraj85"
$ws.Range("H51").Value = "III year CSE B"
$ws.Range("I51").Value = 19679
$ws.Range("J51").Value = "0/4"
$ws.Range("K51").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/788/"
$ws.Range("E52").Value = 48
$ws.Range("F52").Value = "HARISIVAM J J"
$ws.Range("G52").Value = "harisivam_38"
$ws.Range("H52").Value = "III year CSE A"
$ws.Range("I52").Value = 19758
$ws.Range("J52").Value = "0/4"
$ws.Range("K52").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/791/"
$ws.Range("E53").Value = 49
$ws.Range("F53").Value = "KEERTHANA M"
$ws.Range("G53").Value = "user4153n"
$ws.Range("H53").Value = "III year CSE A"
$ws.Range("I53").Value = 19862
$ws.Range("J53").Value = "0/4"
$ws.Range("K53").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/795/"
$ws.Range("E54").Value = 50
$ws.Range("F54").Value = "NAVEEN RAJA S"
$ws.Range("G54").Value = "naveen20cs078"
$ws.Range("H54").Value = "III year CSE B"
$ws.Range("I54").Value = 19918
$ws.Range("J54").Value = "0/4"
$ws.Range("K54").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/797/"
$ws.Range("E55").Value = 51
$ws.Range("F55").Value = "ASHWITHA NOBLE P"
$ws.Range("G55").Value = "ashwitha_noble"
$ws.Range("H55").Value = "III year CSE A"
$ws.Range("I55").Value = 19928
$ws.Range("J55").Value = "0/4"
$ws.Range("K55").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/798/"
$ws.Range("E56").Value = 52
$ws.Range("F56").Value = "ASHISH SINGH"
$ws.Range("G56").Value = "ashish_singh_20cs008"
$ws.Range("H56").Value = "III year CSE A"
$ws.Range("I56").Value = 19971
$ws.Range("J56").Value = "0/4"
$ws.Range("K56").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/799/"
$ws.Range("E57").Value = 53
$ws.Range("F57").Value = "VIDHYA N"
$ws.Range("G57").Value = "vidhya_n"
$ws.Range("H57").Value = "III year CSE C"
$ws.Range("I57").Value = 20029
$ws.Range("J57").Value = "0/4"
$ws.Range("K57").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/802/"
$ws.Range("E58").Value = 54
$ws.Range("F58").Value = "VELLANKI SHAINI"
$ws.Range("G58").Value = "shaini"
$ws.Range("H58").Value = "III year CSE C"
$ws.Range("I58").Value = 20106
$ws.Range("J58").Value = "0/4"
$ws.Range("K58").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/805/"
$ws.Range("E59").Value = 55
$ws.Range("F59").Value = "ABHI NIVESH R"
$ws.Range("G59").Value = "abhi_nivesh_r"
$ws.Range("H59").Value = "III year CSE A"
$ws.Range("I59").Value = 20248
$ws.Range("J59").Value = "0/4"
$ws.Range("K59").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/810/"
$ws.Range("E60").Value = 56
$ws.Range("F60").Value = "FARHEEN A S"
$ws.Range("G60").Value = "farheen_a_s"
$ws.Range("H60").Value = "III year CSE A"
$ws.Range("I60").Value = 20259
$ws.Range("J60").Value = "0/4"
$ws.Range("K60").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/811/"
$ws.Range("E61").Value = 57
$ws.Range("F61").Value = "CHINMAI DEEPIKA M"
$ws.Range("G61").Value = "chinmai_03"
$ws.Range("H61").Value = "III year CSE A"
$ws.Range("I61").Value = 20419
$ws.Range("J61").Value = "0/4"
$ws.Range("K61").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/817/"
$ws.Range("E62").Value = 58
$ws.Range("F62").Value = "MADHUMITHA P"
$ws.Range("G62").Value = "madhu_8523"
$ws.Range("H62").Value = "III year CSE B"
$ws.Range("I62").Value = 20712
$ws.Range("J62").Value = "0/4"
$ws.Range("K62").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/829/"
$ws.Range("E63").Value = 59
$ws.Range("F63").Value = "GAYATHRI M"
$ws.Range("G63").Value = "20cs027"
$ws.Range("H63").Value = "III year CSE A"
$ws.Range("I63").Value = 20754
$ws.Range("J63").Value = "0/4"
$ws.Range("K63").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/831/"
$ws.Range("E64").Value = 60
$ws.Range("F64").Value = "BAGAVATHI ANANDHAN E"
$ws.Range("G64").Value = "anandeswaran"
$ws.Range("H64").Value = "III year CSE A"
$ws.Range("I64").Value = 20801
$ws.Range("J64").Value = "0/4"
$ws.Range("K64").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/833/"
$ws.Range("E65").Value = 61
$ws.Range("F65").Value = "DHANUSHYA R"
$ws.Range("G65").Value = "dhanushyar"
$ws.Range("H65").Value = "III year CSE A"
$ws.Range("I65").Value = 20836
$ws.Range("J65").Value = "0/4"
$ws.Range("K65").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/834/"
$ws.Range("E66").Value = 62
$ws.Range("F66").Value = "SURJITHRAJA A"
$ws.Range("G66").Value = "surjith_07"
$ws.Range("H66").Value = "III year CSE C"
$ws.Range("I66").Value = 20860
$ws.Range("J66").Value = "0/4"
$ws.Range("K66").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/835/"
$ws.Range("E67").Value = 63
$ws.Range("F67").Value = "VARSHINI B"
$ws.Range("G67").Value = "varshini0322"
$ws.Range("H67").Value = "III year CSE C"
$ws.Range("I67").Value = 20968
$ws.Range("J67").Value = "0/4"
$ws.Range("K67").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/839/"
$ws.Range("E68").Value = 64
$ws.Range("F68").Value = "JEEVA JOTHI V M"
$ws.Range("G68").Value = "jeeva_jothi"
$ws.Range("H68").Value = "III year CSE A"
$ws.Range("I68").Value = 21003
$ws.Range("J68").Value = "0/4"
$ws.Range("K68").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/841/"
$ws.Range("E69").Value = 65
$ws.Range("F69").Value = "ARTHI V"
$ws.Range("G69").Value = "20cs006"
$ws.Range("H69").Value = "III year CSE A"
$ws.Range("I69").Value = 21074
$ws.Range("J69").Value = "0/4"
$ws.Range("K69").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/843/"
$ws.Range("E70").Value = 66
$ws.Range("F70").Value = "JAGAPREETHA B"
$ws.Range("G70").Value = "jagapreetha_b"
$ws.Range("H70").Value = "III year CSE A"
$ws.Range("I70").Value = 21132
$ws.Range("J70").Value = "0/4"
$ws.Range("K70").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/846/"
$ws.Range("E71").Value = 67
$ws.Range("F71").Value = "AKSHAY KRISHNA N"
$ws.Range("G71").Value = "akshay_krishna"
$ws.Range("H71").Value = "III year CSE A"
$ws.Range("I71").Value = 21156
$ws.Range("J71").Value = "0/4"
$ws.Range("K71").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/847/"
$ws.Range("E72").Value = 68
$ws.Range("F72").Value = "SUGANTHI M"
$ws.Range("G72").Value = "suganthi_m"
$ws.Range("H72").Value = "III year CSE C"
$ws.Range("I72").Value = 21182
$ws.Range("J72").Value = "0/4"
$ws.Range("K72").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/848/"
$ws.Range("E73").Value = 69
$ws.Range("F73").Value = "SUBAHARINI"
$ws.Range("G73").Value = "subhaharini_s"
$ws.Range("H73").Value = "III year CSE C"
$ws.Range("I73").Value = 21194
$ws.Range("J73").Value = "0/4"
$ws.Range("K73").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/848/"
$ws.Range("E74").Value = 70
$ws.Range("F74").Value = "SUBANANTHITHA K"
$ws.Range("G74").Value = "suba17"
$ws.Range("H74").Value = "III year CSE C"
$ws.Range("I74").Value = 21204
$ws.Range("J74").Value = "0/4"
$ws.Range("K74").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/849/"
$ws.Range("E75").Value = 71
$ws.Range("F75").Value = "SAKTHIMAHESWARAN U S"
$ws.Range("G75").Value = "sakthimaheswaran_2002"
$ws.Range("H75").Value = "III year CSE B"
$ws.Range("I75").Value = 21238
$ws.Range("J75").Value = "0/4"
$ws.Range("K75").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/850/"
$ws.Range("E76").Value = 72
$ws.Range("F76").Value = "MARIUSH RUFIN P"
$ws.Range("G76").Value = "mariush123"
$ws.Range("H76").Value = "III year CSE B"
$ws.Range("I76").Value = 21260
$ws.Range("J76").Value = "0/4"
$ws.Range("K76").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/851/"
$ws.Range("E77").Value = 73
$ws.Range("F77").Value = "GNANA JOTHI T"
$ws.Range("G77").Value = "gnanajothi"
$ws.Range("H77").Value = "III year CSE A"
$ws.Range("I77").Value = 21287
$ws.Range("J77").Value = "0/4"
$ws.Range("K77").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/852/"
$ws.Range("E78").Value = 74
$ws.Range("F78").Value = "GOLUGURI YASWANTH DURGA REDDY"
$ws.Range("G78").Value = "yaswanth1825"
$ws.Range("H78").Value = "III year CSE A"
$ws.Range("I78").Value = 21294
$ws.Range("J78").Value = "0/4"
$ws.Range("K78").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/852/"
$ws.Range("E79").Value = 75
$ws.Range("F79").Value = "MOHAN RAAM I P"
$ws.Range("G79").Value = "mohanraam21"
$ws.Range("H79").Value = "III year CSE B"
$ws.Range("I79").Value = 21300
$ws.Range("J79").Value = "0/4"
$ws.Range("K79").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/852/"
$ws.Range("E80").Value = 76
$ws.Range("F80").Value = "GOVINDHA RAMANATHAN S"
$ws.Range("G80").Value = "user5780js"
$ws.Range("H80").Value = "III year CSE A"
$ws.Range("I80").Value = 21305
$ws.Range("J80").Value = "0/4"
$ws.Range("K80").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/853/"
$ws.Range("E81").Value = 77
$ws.Range("F81").Value = "NAGARJUNAN G"
$ws.Range("G81").Value = "arj742"
$ws.Range("H81").Value = "III year CSE B"
$ws.Range("I81").Value = 21376
$ws.Range("J81").Value = "0/4"
$ws.Range("K81").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/856/"
$ws.Range("E82").Value = 78
$ws.Range("F82").Value = "PREETHIKA P"
$ws.Range("G82").Value = "preethika_p"
$ws.Range("H82").Value = "III year CSE B"
$ws.Range("I82").Value = 21417
$ws.Range("J82").Value = "0/4"
$ws.Range("K82").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/857/"
$ws.Range("E83").Value = 79
$ws.Range("F83").Value = "SHANMUGA PRIYA R"
$ws.Range("G83").Value = "shanmugapriya126"
$ws.Range("H83").Value = "III year CSE C"
$ws.Range("I83").Value = 21451
$ws.Range("J83").Value = "0/4"
$ws.Range("K83").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/859/"
$ws.Range("E84").Value = 80
$ws.Range("F84").Value = "SHOBANA M"
$ws.Range("G84").Value = "shobana03"
$ws.Range("H84").Value = "III year CSE C"
$ws.Range("I84").Value = 21458
$ws.Range("J84").Value = "0/4"
$ws.Range("K84").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/859/"
$ws.Range("E85").Value = 81
$ws.Range("F85").Value = "BLESSED JEBERSON J Q"
$ws.Range("G85").Value = "blessed_jeberson"
$ws.Range("H85").Value = "III year CSE A"
$ws.Range("I85").Value = 21482
$ws.Range("J85").Value = "0/4"
$ws.Range("K85").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/860/"
$ws.Range("E86").Value = 82
$ws.Range("F86").Value = "BALAJIRAM N"
$ws.Range("G86").Value = "user2849ec"
$ws.Range("H86").Value = "III year CSE A"
$ws.Range("I86").Value = 21486
$ws.Range("J86").Value = "0/4"
$ws.Range("K86").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/860/"
$ws.Range("E87").Value = 83
$ws.Range("F87").Value = "NIKHILESH S"
$ws.Range("G87").Value = "nikhilnikhil"
$ws.Range("H87").Value = "III year CSE B"
$ws.Range("I87").Value = 21493
$ws.Range("J87").Value = "0/4"
$ws.Range("K87").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/860/"
$ws.Range("E88").Value = 84
$ws.Range("F88").Value = "RANJITHKUMAR D"
$ws.Range("G88").Value = "ranjithkumar7"
$ws.Range("H88").Value = "III year CSE B"
$ws.Range("I88").Value = 21501
$ws.Range("J88").Value = "0/4"
$ws.Range("K88").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/861/"
$ws.Range("E89").Value = 85
$ws.Range("F89").Value = "NIHAR PARVEEN A"
$ws.Range("G89").Value = "nihar_2003"
$ws.Range("H89").Value = "III year CSE B"
$ws.Range("I89").Value = 21503
$ws.Range("J89").Value = "0/4"
$ws.Range("K89").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/861/"
$ws.Range("E90").Value = 86
$ws.Range("F90").Value = "SANJEEVA KUMAR M"
$ws.Range("G90").Value = "sanjeevakumar"
$ws.Range("H90").Value = "III year CSE C"
$ws.Range("I90").Value = 21623
$ws.Range("J90").Value = "0/4"
$ws.Range("K90").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/865/"
$ws.Range("E91").Value = 87
$ws.Range("F91").Value = "KIRTHI VIGNESH G"
$ws.Range("G91").Value = "kirthi_vignesh_g"
$ws.Range("H91").Value = "III year CSE B"
$ws.Range("I91").Value = 21642
$ws.Range("J91").Value = "0/4"
$ws.Range("K91").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/866/"
$ws.Range("E92").Value = 88
$ws.Range("F92").Value = "SUGAVANAESH S"
$ws.Range("G92").Value = "sugavanaesh"
$ws.Range("H92").Value = "III year CSE C"
$ws.Range("I92").Value = 21658
$ws.Range("J92").Value = "0/4"
$ws.Range("K92").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/867/"
$ws.Range("E93").Value = 89
$ws.Range("F93").Value = "ILAMUGUNTHAN N"
$ws.Range("G93").Value = "ilamugunthan"
$ws.Range("H93").Value = "III year CSE A"
$ws.Range("I93").Value = 21660
$ws.Range("J93").Value = "0/4"
$ws.Range("K93").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/867/"
$ws.Range("E94").Value = 90
$ws.Range("F94").Value = "DHARANI P"
$ws.Range("G94").Value = "20cs021"
$ws.Range("H94").Value = "III year CSE A"
$ws.Range("I94").Value = 21664
$ws.Range("J94").Value = "0/4"
$ws.Range("K94").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/867/"
$ws.Range("E95").Value = 91
$ws.Range("F95").Value = "PRASANNA K M"
$ws.Range("G95").Value = "_prasanna_k_m_"
$ws.Range("H95").Value = "III year CSE B"
$ws.Range("I95").Value = 21667
$ws.Range("J95").Value = "0/4"
$ws.Range("K95").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/867/"
$ws.Range("E96").Value = 92
$ws.Range("F96").Value = "PRANESH S"
$ws.Range("G96").Value = "pranesh_22"
$ws.Range("H96").Value = "III year CSE B"
$ws.Range("I96").Value = 21670
$ws.Range("J96").Value = "0/4"
$ws.Range("K96").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/867/"
$ws.Range("E97").Value = 93
$ws.Range("F97").Value = "SIVANARAYAN J"
$ws.Range("G97").Value = "sivanarayan"
$ws.Range("H97").Value = "III year CSE C"
$ws.Range("I97").Value = 21674
$ws.Range("J97").Value = "0/4"
$ws.Range("K97").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/867/"
$ws.Range("E98").Value = 94
$ws.Range("F98").Value = "VISHAL KHUMAR R"
$ws.Range("G98").Value = "vishalkhumar_r"
$ws.Range("H98").Value = "III year CSE C"
$ws.Range("I98").Value = 21675
$ws.Range("J98").Value = "0/4"
$ws.Range("K98").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/867/"
$ws.Range("E99").Value = 95
$ws.Range("F99").Value = "VISHNUPRABHU R"
$ws.Range("G99").Value = "user1652be"
$ws.Range("H99").Value = "III year CSE C"
$ws.Range("I99").Value = 21686
$ws.Range("J99").Value = "0/4"
$ws.Range("K99").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/868/"
$ws.Range("E100").Value = 96
$ws.Range("F100").Value = "SHAVAL KHAN M"
$ws.Range("G100").Value = "shavalkhan"
$ws.Range("H100").Value = "III year CSE C"
$ws.Range("I100").Value = 21741
$ws.Range("J100").Value = "0/4"
$ws.Range("K100").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/870/"
$ws.Range("E101").Value = 97
$ws.Range("F101").Value = "KARTHEKEIAN K R"
$ws.Range("G101").Value = "karthekeian"
$ws.Range("H101").Value = "III year CSE A"
$ws.Range("I101").Value = 21761
$ws.Range("J101").Value = "0/4"
$ws.Range("K101").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/871/"
$ws.Range("E102").Value = 98
$ws.Range("F102").Value = "DIVYA PRIYA S"
$ws.Range("G102").Value = "divyapriya_s"
$ws.Range("H102").Value = "III year CSE A"
$ws.Range("I102").Value = 21765
$ws.Range("J102").Value = "0/4"
$ws.Range("K102").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/871/"
$ws.Range("E103").Value = 99
$ws.Range("F103").Value = "SAI SIDDHARTH"
$ws.Range("G103").Value = "sai_sid_2002"
$ws.Range("H103").Value = "III year CSE B"
$ws.Range("I103").Value = 21766
$ws.Range("J103").Value = "0/4"
$ws.Range("K103").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/871/"
$ws.Range("E104").Value = 100
$ws.Range("F104").Value = "GEETHANJALI G"
$ws.Range("G104").Value = "geethanjali_g"
$ws.Range("H104").Value = "III year CSE A"
$ws.Range("I104").Value = 21771
$ws.Range("J104").Value = "0/4"
$ws.Range("K104").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/871/"
$ws.Range("E105").Value = 101
$ws.Range("F105").Value = "SELVENDHIRAN R"
$ws.Range("G105").Value = "selva_7"
$ws.Range("H105").Value = "III year CSE C"
$ws.Range("I105").Value = 21782
$ws.Range("J105").Value = "0/4"
$ws.Range("K105").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/872/"
$ws.Range("E106").Value = 102
$ws.Range("F106").Value = "KISHOREVEL I V"
$ws.Range("G106").Value = "vel_02"
$ws.Range("H106").Value = "III year CSE B"
$ws.Range("I106").Value = 21811
$ws.Range("J106").Value = "0/4"
$ws.Range("K106").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/873/"
$ws.Range("E107").Value = 103
$ws.Range("F107").Value = "GOWSHIKAN S"
$ws.Range("G107").Value = "gowshikan_s"
$ws.Range("H107").Value = "III year CSE A"
$ws.Range("I107").Value = 21822
$ws.Range("J107").Value = "0/4"
$ws.Range("K107").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/873/"
$ws.Range("E108").Value = 104
$ws.Range("F108").Value = "SREE GOUSHIK RAJAA R"
$ws.Range("G108").Value = "bloodlust"
$ws.Range("H108").Value = "III year CSE C"
$ws.Range("I108").Value = 21842
$ws.Range("J108").Value = "0/4"
$ws.Range("K108").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/874/"
$ws.Range("E109").Value = 105
$ws.Range("F109").Value = "GOKUL S"
$ws.Range("G109").Value = "gokulcsr"
$ws.Range("H109").Value = "III year CSE A"
$ws.Range("I109").Value = 21850
$ws.Range("J109").Value = "0/4"
$ws.Range("K109").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/874/"
$ws.Range("E110").Value = 106
$ws.Range("F110").Value = "DHANUSH B"
$ws.Range("G110").Value = "user4297oq"
$ws.Range("H110").Value = "III year CSE A"
$ws.Range("I110").Value = 21857
$ws.Range("J110").Value = "0/4"
$ws.Range("K110").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/875/"
$ws.Range("E111").Value = 107
$ws.Range("F111").Value = "KARTHIK V"
$ws.Range("G111").Value = "karthikv27"
$ws.Range("H111").Value = "III year CSE A"
$ws.Range("I111").Value = 21863
$ws.Range("J111").Value = "0/4"
$ws.Range("K111").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/875/"
$ws.Range("E112").Value = 108
$ws.Range("F112").Value = "SRIMATHI G"
$ws.Range("G112").Value = "srimathi140"
$ws.Range("H112").Value = "III year CSE C"
$ws.Range("I112").Value = 21868
$ws.Range("J112").Value = "0/4"
$ws.Range("K112").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/875/"
$ws.Range("E113").Value = 109
$ws.Range("F113").Value = "CHARAN SRI KRISHNA S"
$ws.Range("G113").Value = "charansri795"
$ws.Range("H113").Value = "III year CSE A"
$ws.Range("I113").Value = 21873
$ws.Range("J113").Value = "0/4"
$ws.Range("K113").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/875/"
$ws.Range("E114").Value = 110
$ws.Range("F114").Value = "SARVESH V V"
$ws.Range("G114").Value = "sarvesh01"
$ws.Range("H114").Value = "III year CSE C"
$ws.Range("I114").Value = 21881
$ws.Range("J114").Value = "0/4"
$ws.Range("K114").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/876/"
$ws.Range("E115").Value = 111
$ws.Range("F115").Value = "KAVINKUMAR S"
$ws.Range("G115").Value = "kavinkumar_s"
$ws.Range("H115").Value = "III year CSE A"
$ws.Range("I115").Value = 21883
$ws.Range("J115").Value = "0/4"
$ws.Range("K115").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/876/"
$ws.Range("E116").Value = 112
$ws.Range("F116").Value = "VINISHA M"
$ws.Range("G116").Value = "vinisha_m"
$ws.Range("H116").Value = "III year CSE C"
$ws.Range("I116").Value = 21894
$ws.Range("J116").Value = "0/4"
$ws.Range("K116").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/876/"
$ws.Range("E117").Value = 113
$ws.Range("F117").Value = "MOHANAPRASATH M"
$ws.Range("G117").Value = "mohan_09"
$ws.Range("H117").Value = "III year CSE B"
$ws.Range("I117").Value = 21929
$ws.Range("J117").Value = "0/4"
$ws.Range("K117").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/878/"
$ws.Range("E118").Value = 114
$ws.Range("F118").Value = "SREESHARAN N"
$ws.Range("G118").Value = "sharan20cs138"
$ws.Range("H118").Value = "III year CSE C"
$ws.Range("I118").Value = 21987
$ws.Range("J118").Value = "0/4"
$ws.Range("K118").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/880/"
$ws.Range("E119").Value = 115
$ws.Range("F119").Value = "SANJAY N"
$ws.Range("G119").Value = "sanjay_n_2003"
$ws.Range("H119").Value = "III year CSE B"
$ws.Range("I119").Value = 22158
$ws.Range("J119").Value = "0/4"
$ws.Range("K119").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/887/"
$ws.Range("E120").Value = 116
$ws.Range("F120").Value = "KEERTHANA P"
$ws.Range("G120").Value = "user5187jk"
$ws.Range("H120").Value = "III year CSE A"
$ws.Range("I120").Value = 22182
$ws.Range("J120").Value = "0/4"
$ws.Range("K120").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/888/"
$ws.Range("E121").Value = 117
$ws.Range("F121").Value = "ABBINAVU T"
$ws.Range("G121").Value = "20cs001"
$ws.Range("H121").Value = "III year CSE A"
$ws.Range("I121").Value = 22389
$ws.Range("J121").Value = "0/4"
$ws.Range("K121").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/896/"
$ws.Range("E122").Value = 118
$ws.Range("F122").Value = "GOWREESH A M"
$ws.Range("G122").Value = "gowreesh_a_m"
$ws.Range("H122").Value = "III year CSE A"
$ws.Range("I122").Value = 22413
$ws.Range("J122").Value = "0/4"
$ws.Range("K122").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/897/"
$ws.Range("E123").Value = 119
$ws.Range("F123").Value = "DINESH P"
$ws.Range("G123").Value = "dineshpraba"
$ws.Range("H123").Value = "III year CSE A"
$ws.Range("I123").Value = 22503
$ws.Range("J123").Value = "0/4"
$ws.Range("K123").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/901/"
$ws.Range("E124").Value = 120
$ws.Range("F124").Value = "HEMRAJ KUMAR V S"
$ws.Range("G124").Value = "hems348"
$ws.Range("H124").Value = "III year CSE A"
$ws.Range("I124").Value = 22507
$ws.Range("J124").Value = "0/4"
$ws.Range("K124").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/901/"
$ws.Range("E125").Value = 121
$ws.Range("F125").Value = "RAHUL RAJ R"
$ws.Range("G125").Value = "rahulrajr"
$ws.Range("H125").Value = "III year CSE B"
$ws.Range("I125").Value = 22559
$ws.Range("J125").Value = "0/4"
$ws.Range("K125").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/903/"
$ws.Range("E126").Value = 122
$ws.Range("F126").Value = "JAYANTH N J"
$ws.Range("G126").Value = "jayanth_n_j"
$ws.Range("H126").Value = "III year CSE A"
$ws.Range("I126").Value = 22611
$ws.Range("J126").Value = "0/4"
$ws.Range("K126").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/905/"
$ws.Range("E127").Value = 123
$ws.Range("F127").Value = "ESHA MALAVIKA V S"
$ws.Range("G127").Value = "eshamalavika"
$ws.Range("H127").Value = "III year CSE A"
$ws.Range("I127").Value = 22626
$ws.Range("J127").Value = "0/4"
$ws.Range("K127").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/906/"
$ws.Range("E128").Value = 124
$ws.Range("F128").Value = "VINUDHARSHINI R"
$ws.Range("G128").Value = "vinudharshini_r"
$ws.Range("H128").Value = "III year CSE C"
$ws.Range("I128").Value = 22632
$ws.Range("J128").Value = "0/4"
$ws.Range("K128").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/906/"
$ws.Range("E129").Value = 125
$ws.Range("F129").Value = "VARSHIKHA N R"
$ws.Range("G129").Value = "varshikha_nr"
$ws.Range("H129").Value = "III year CSE C"
$ws.Range("I129").Value = 22654
$ws.Range("J129").Value = "0/4"
$ws.Range("K129").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/907/"
$ws.Range("E130").Value = 126
$ws.Range("F130").Value = "KANIPRIYA R"
$ws.Range("G130").Value = "kani_03"
$ws.Range("H130").Value = "III year CSE A"
$ws.Range("I130").Value = 22702
$ws.Range("J130").Value = "0/4"
$ws.Range("K130").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/909/"
$ws.Range("E131").Value = 127
$ws.Range("F131").Value = "DIVAKAR V S"
$ws.Range("G131").Value = "20cs023"
$ws.Range("H131").Value = "III year CSE A"
$ws.Range("I131").Value = 22711
$ws.Range("J131").Value = "0/4"
$ws.Range("K131").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/909/"
$ws.Range("E132").Value = 128
$ws.Range("F132").Value = "ARUN KUMAR P"
$ws.Range("G132").Value = "arun_kumar_p"
$ws.Range("H132").Value = "III year CSE A"
$ws.Range("I132").Value = 22720
$ws.Range("J132").Value = "0/4"
$ws.Range("K132").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/909/"

# ===================== Sheet: III year it =====================
$ws = $wb.Worksheets.Item("III year it")
$ws.Range("I5").Value = 13732
$ws.Range("K5").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/550/"
$ws.Range("I6").Value = 19940
$ws.Range("K6").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/798/"
$ws.Range("I7").Value = 21631
$ws.Range("K7").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/866/"

# ===================== Sheet: III year csbs =====================
$ws = $wb.Worksheets.Item("III year csbs")
$ws.Range("I5").Value = 19557
$ws.Range("K5").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/783/"
$ws.Range("I6").Value = 20818
$ws.Range("K6").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/833/"
$ws.Range("I7").Value = 20822
$ws.Range("K7").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/833/"
$ws.Range("I8").Value = 22350
$ws.Range("K8").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/894/"

# ===================== Sheet: III year ai & ds =====================
$ws = $wb.Worksheets.Item("III year ai & ds")
$ws.Range("I5").Value = 19457
$ws.Range("K5").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/779/"
$ws.Range("I6").Value = 19464
$ws.Range("K6").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/779/"
$ws.Range("I7").Value = 19976
$ws.Range("K7").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/800/"
$ws.Range("I8").Value = 20755
$ws.Range("K8").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/831/"
$ws.Range("I9").Value = 20776
$ws.Range("K9").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/832/"
$ws.Range("I10").Value = 20804
$ws.Range("K10").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/833/"
$ws.Range("I11").Value = 21164
$ws.Range("K11").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/847/"
$ws.Range("I12").Value = 21292
$ws.Range("K12").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/852/"
$ws.Range("I13").Value = 21323
$ws.Range("K13").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/853/"
$ws.Range("I14").Value = 21469
$ws.Range("K14").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/859/"
$ws.Range("I15").Value = 21515
$ws.Range("K15").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/861/"
$ws.Range("I16").Value = 21628
$ws.Range("K16").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/866/"
$ws.Range("I17").Value = 21654
$ws.Range("K17").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/867/"
$ws.Range("I18").Value = 21656
$ws.Range("K18").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/867/"
$ws.Range("I19").Value = 21663
$ws.Range("K19").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/867/"
$ws.Range("I20").Value = 21700
$ws.Range("K20").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/868/"
$ws.Range("I21").Value = 21745
$ws.Range("K21").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/870/"
$ws.Range("I22").Value = 21770
$ws.Range("K22").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/871/"
$ws.Range("I23").Value = 21871
$ws.Range("K23").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/875/"

# ===================== Sheet: III year ece =====================
$ws = $wb.Worksheets.Item("III year ece")
$ws.Range("I5").Value = 14362
$ws.Range("K5").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/575/"
$ws.Range("I6").Value = 15320
$ws.Range("K6").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/613/"
$ws.Range("I7").Value = 18436
$ws.Range("K7").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/738/"
$ws.Range("I8").Value = 19650
$ws.Range("K8").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/786/"
$ws.Range("I9").Value = 19734
$ws.Range("K9").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/790/"
$ws.Range("I10").Value = 21424
$ws.Range("K10").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/857/"
$ws.Range("I11").Value = 21746
$ws.Range("K11").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/870/"

# ===================== Sheet: III year cce =====================
$ws = $wb.Worksheets.Item("III year cce")
$ws.Range("I5").Value = 17710
$ws.Range("K5").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/709/"
$ws.Range("I6").Value = 21046
$ws.Range("K6").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/842/"
$ws.Range("I7").Value = 21169
$ws.Range("K7").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/847/"
$ws.Range("I8").Value = 21253
$ws.Range("K8").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/851/"
$ws.Range("I9").Value = 21271
$ws.Range("K9").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/851/"
$ws.Range("I10").Value = 21436
$ws.Range("K10").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/858/"
$ws.Range("I11").Value = 21462
$ws.Range("K11").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/859/"
$ws.Range("I12").Value = 21599
$ws.Range("K12").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/864/"
$ws.Range("I13").Value = 21617
$ws.Range("K13").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/865/"
$ws.Range("I14").Value = 21627
$ws.Range("K14").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/866/"
$ws.Range("I15").Value = 21640
$ws.Range("K15").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/866/"
$ws.Range("I16").Value = 21647
$ws.Range("K16").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/866/"
$ws.Range("I17").Value = 22429
$ws.Range("K17").Value = "https://leetcode.com/contest/weekly-contest-308/ranking/897/"

